{"js": "const replacements = [\n  [\"which were among\", \"which are among\"],\n  [\"influenced by the new approach\", \"influenced by this novel approach\"],\n  [\"My doctoral research process has begun\", \"My doctoral research journey has begun\"],\n  [\"found their way into an electronic voting academic proposal. All except\",\n   \"found their way into electronic voting academic proposals. All except\"],\n  [\"privacy-preserving method to commit data to a blockchain\",\n   \"privacy-preserving method for commiting information to a blockchain\"],\n  [\"decentralised e-voting proposal. Yet\", \"decentralised e-voting platform. Yet\"],\n  [\"using them to develop a NFT-based e-voting proposal\",\n   \"using them to develop an NFT-based e-voting proposal\"],\n  [\"establishes a period of 6 months minimum to a maximum of one year where\",\n   \"establishes a period (from a minimum of 6 months to a maximum of one year) where\"],\n  [\"I originally met some faculty from the University of Surrey during a prior research collaboration.\",\n   \"I originally met some faculty members from the University of Surrey during a prior research collaboration on the \u201cTransition Guardian\u201d line of research.\"],\n  [\"a renowed centre\", \"a renowned centre\"],\n  [\"I\u2019m planning to have a NFT architectural\", \"I\u2019m planning to have an NFT architectural\"],\n  [\"development of a NFT-based e-voting system\", \"development of an NFT-based e-voting system\"],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + findText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Tables ready, text needs work\n# Apply the wording revisions to the Research Plan narrative paragraphs.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"No match found for: $findText\"\n    }\n}\n\n# Paragraph 1 (research summary)\nReplace-Text \"which were among\" \"which are among\"\nReplace-Text \"influenced by the new approach\" \"influenced by this novel approach\"\nReplace-Text \"My doctoral research process has begun\" \"My doctoral research journey has begun\"\nReplace-Text \"found their way into an electronic voting academic proposal. All except\" \"found their way into electronic voting academic proposals. All except\"\nReplace-Text \"privacy-preserving method to commit data to a blockchain\" \"privacy-preserving method for commiting information to a blockchain\"\nReplace-Text \"decentralised e-voting proposal. Yet\" \"decentralised e-voting platform. Yet\"\nReplace-Text \"using them to develop a NFT-based e-voting proposal\" \"using them to develop an NFT-based e-voting proposal\"\n\n# Paragraph 2 (Surrey collaboration)\nReplace-Text \"establishes a period of 6 months minimum to a maximum of one year where\" \"establishes a period (from a minimum of 6 months to a maximum of one year) where\"\nReplace-Text \"I originally met some faculty from the University of Surrey during a prior research collaboration.\" \"I originally met some faculty members from the University of Surrey during a prior research collaboration on the \u201cTransition Guardian\u201d line of research.\"\nReplace-Text \"a renowed centre\" \"a renowned centre\"\nReplace-Text \"I\u2019m planning to have a NFT architectural\" \"I\u2019m planning to have an NFT architectural\"\nReplace-Text \"development of a NFT-based e-voting system\" \"development of an NFT-based e-voting system\"\n"}
